$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update C8:C13 with new step labels (values for rows 2-7 remain unchanged)
$ws.Range("C8").Value = "F.0.1.21_3"
$ws.Range("C9").Value = "F.0.1.21_4"
$ws.Range("C10").Value = "F.0.1.22_1"
$ws.Range("C11").Value = "F.0.1.22_2"
$ws.Range("C12").Value = "F.0.1.22_3"
$ws.Range("C13").Value = "F.0.1.22_4"

# Update H2:H13 from 1 to 2
$ws.Range("H2:H13").Value = 2

# Update the active selection to J7
$ws.Range("J7").Select()
